$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix typo "serorreversion" -> "seroreversion" in the second table header (row 8)
$ws.Range("F8").Value = "Tasa de seroreversion"
$ws.Range("G8").Value = "Tasa de seroreversión Rhat"

# Update the selection to span the full second table (A8:G12)
$ws.Range("A8:G12").Select()
